$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text storage so that
# values such as "1.000" or "7.139" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.360.49'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.847.11'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").Value = '240.35'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '0.6285'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.07573'
$ws.Range("D9").Value = '0.2910'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").Value = '24.60'
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("D11").Value = '0.07752'
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").Value = '1.845.89'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '5.015'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '0.6791'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = '0.00001045'
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("D16").Value = '83.07'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '6.109'
$ws.Range("E17").Value = '  -1.01%  '
$ws.Range("D18").Value = '29.322.95'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = '228.93'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '12.34'
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '7.420'
$ws.Range("E22").Value = '  -0.37%  '
$ws.Range("D23").Value = '0.9995'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").Value = '159.01'
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("D25").Value = '0.1389'
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = '8.426'
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").Value = '17.66'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = '1.424'
$ws.Range("E28").Value = '  +7.49%  '
$ws.Range("D29").Value = '1.471'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = '0.05671'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = '4.111'
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").Value = '4.037'
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '1.824'
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.155'
$ws.Range("E34").Value = '  -0.86%  '
$ws.Range("D35").Value = '0.6956'
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("D36").Value = '2.581'
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("D37").Value = '0.01827'
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").Value = '1.234.94'
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").Value = '2.723'
$ws.Range("E39").Value = '  -2.25%  '
$ws.Range("D40").Value = '6.375'
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("D41").Value = '0.8968'
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").Value = '0.9999'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '101.27'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '65.44'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '7.139'
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '0.4001'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = '0.1152'
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").Value = '8.962'
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").Value = '1.672'
$ws.Range("E50").Value = '  -0.58%  '
$ws.Range("D51").Value = '0.05697'
$ws.Range("E51").Value = '  -0.38%  '
